$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the stray empty B84 cell (it currently holds an empty inline string).
$ws.Range("B84").ClearContents()

# Append the four new newsbot rows (85-88).
# Row 85
$ws.Range("A85").Value = '05/01/2026 12:49:52'
$ws.Range("B85").Value = '05/01 12:39'
$ws.Range("C85").Value = 'g1 > Política'
$ws.Range("D85").Value = 'Sem citar caso Master, entidades do setor financeiro divulgam nota em defesa do BC: ''Plena confiança'''
$ws.Range("E85").Value = 'https://g1.globo.com/economia/noticia/2026/01/05/sem-citar-caso-master-entidades-do-setor-financeiro-divulgam-nota-em-defesa-do-bc-plena-confianca.ghtml'
$ws.Range("F85").Value = 'tcu'
$ws.Range("G85").Value = ' autoridade monetária. 
Neste início de ano, o presidente do Tribunal de Contas da União (TCU), ministro Vital do Rêgo Filho, determinou uma inspeção técnica em documentos referentes '
# The embedded line break above makes Excel auto-set a custom row height;
# restore the sheet's default (no explicit row-height override), matching the source.
$ws.Rows.Item(85).AutoFit()

# Row 86
$ws.Range("A86").Value = '05/01/2026 12:49:53'
$ws.Range("B86").Value = '05/01 12:28'
$ws.Range("C86").Value = 'g1 > Política'
$ws.Range("D86").Value = 'TCU confirma autorização de inspeção no Banco Central sobre liquidação do Master'
$ws.Range("E86").Value = 'https://g1.globo.com/economia/noticia/2026/01/05/tcu-confirma-autorizacao-de-inspecao-no-bc-sobre-liquidacao-do-master.ghtml'
$ws.Range("F86").Value = 'câmara'
$ws.Range("G86").Value = 'than de Jesus acolheu um pedido do Ministério Público do TCU e da liderança da minoria na Câmara dos Deputados e pediu esclarecimentos ao Banco Central sobre os motivos da liquidação do '

# Row 87
$ws.Range("A87").Value = '05/01/2026 12:49:54'
$ws.Range("B87").Value = '05/01 12:26'
$ws.Range("C87").Value = 'Folha de S.Paulo - Mercado - Principal'
$ws.Range("D87").Value = 'Relator do Master no TCU diz manter no radar cautelar em decisão sobre inspeção no BC'
$ws.Range("E87").Value = 'https://redir.folha.com.br/redir/online/mercado/rss091/*https://www1.folha.uol.com.br/mercado/2026/01/relator-do-master-no-tcu-diz-manter-no-radar-cautelar-em-decisao-sobre-inspecao-no-bc.shtml'
$ws.Range("F87").Value = 'tcu'
$ws.Range("G87").Value = 'O relator do caso Master no &lt;a href="https://www1.folha.uol.com.br/folha-topicos/&lt;b&gt;tcu&lt;/b&gt;/"&gt;TCU&lt;/a&gt; (Tribunal de Contas da União), ministro Jhonatan de Jesus, determinou nesta seg'

# Row 88
$ws.Range("A88").Value = '05/01/2026 12:49:55'
$ws.Range("B88").Value = '05/01 10:37'
$ws.Range("C88").Value = 'Folha de S.Paulo - Mercado - Principal'
$ws.Range("D88").Value = 'Governo tem maior crescimento no número de servidores em dez anos, com 19 mil funcionários a mais'
$ws.Range("E88").Value = 'https://redir.folha.com.br/redir/online/mercado/rss091/*https://www1.folha.uol.com.br/mercado/2026/01/governo-tem-maior-crescimento-no-numero-de-servidores-em-dez-anos-com-19-mil-funcionarios-a-mais.shtml'
$ws.Range("F88").Value = 'lula'
$ws.Range("G88").Value = 'O governo do presidente Luiz Inácio &lt;a href="https://www1.folha.uol.com.br/folha-topicos/&lt;b&gt;lula&lt;/b&gt;/" rel="" target=""&gt;Lula &lt;/a&gt;da Silva (&lt;a href="https://www1.folha.uol.com.br/folha-topico'

